# Auto-generated script to apply meteocat data update (2026-02-17 22:20 run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-17 22:18:36"
$ws.Range("E3").Value = "2026-02-17 22:18:38"
$ws.Range("E4").Value = "2026-02-17 22:18:41"
$ws.Range("J4").Value = "1018.7 hPa"
$ws.Range("E5").Value = "2026-02-17 22:18:43"
$ws.Range("E6").Value = "2026-02-17 22:18:46"
$ws.Range("J6").Value = "1018.6 hPa"
$ws.Range("E7").Value = "2026-02-17 22:18:48"
$ws.Range("J7").Value = "1018.5 hPa"
$ws.Range("E8").Value = "2026-02-17 22:18:51"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "76%"
$ws.Range("I8").Copy()
$ws.Range("H8").PasteSpecial(-4122)
$ws.Range("J8").Value = "1018.4 hPa"
$ws.Range("E9").Value = "2026-02-17 22:18:53"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "61%"
$ws.Range("I9").Copy()
$ws.Range("H9").PasteSpecial(-4122)
$ws.Range("N9").Value = "7.1 °C 21:37 TU"
$ws.Range("O9").Value = "12.1 °C"
$ws.Range("E10").Value = "2026-02-17 22:18:56"
$ws.Range("E11").Value = "2026-02-17 22:18:58"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "54%"
$ws.Range("I11").Copy()
$ws.Range("H11").PasteSpecial(-4122)
$ws.Range("O11").Value = "7.2 °C"
$ws.Range("E12").Value = "2026-02-17 22:19:01"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "63%"
$ws.Range("I12").Copy()
$ws.Range("H12").PasteSpecial(-4122)
$ws.Range("N12").Value = "8.6 °C 21:59 TU"
$ws.Range("O12").Value = "12.4 °C"
$ws.Range("E13").Value = "2026-02-17 22:19:03"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "47%"
$ws.Range("I13").Copy()
$ws.Range("H13").PasteSpecial(-4122)
$ws.Range("E14").Value = "2026-02-17 22:19:05"
$ws.Range("E15").Value = "2026-02-17 22:19:08"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "60%"
$ws.Range("I15").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("O15").Value = "11.8 °C"
$ws.Range("E16").Value = "2026-02-17 22:19:10"
$ws.Range("M16").Value = "0.2 °C 21:59 TU"
$ws.Range("O16").Value = "-3.2 °C"
$ws.Range("E17").Value = "2026-02-17 22:19:13"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "83%"
$ws.Range("I17").Copy()
$ws.Range("H17").PasteSpecial(-4122)
$ws.Range("E18").Value = "2026-02-17 22:19:15"
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").Value = "80%"
$ws.Range("I18").Copy()
$ws.Range("H18").PasteSpecial(-4122)
$ws.Range("O18").Value = "10.2 °C"
$ws.Range("E19").Value = "2026-02-17 22:19:18"
$ws.Range("E20").Value = "2026-02-17 22:19:20"
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = "67%"
$ws.Range("I20").Copy()
$ws.Range("H20").PasteSpecial(-4122)
$ws.Range("E21").Value = "2026-02-17 22:19:23"
$ws.Range("J21").Value = "1017.2 hPa"
$ws.Range("O21").Value = "9.5 °C"
$ws.Range("E22").Value = "2026-02-17 22:19:25"
$ws.Range("E23").Value = "2026-02-17 22:19:27"
$ws.Range("M23").Value = "-0.2 °C 21:53 TU"
$ws.Range("O23").Value = "-3.6 °C"
$ws.Range("E24").Value = "2026-02-17 22:19:30"
$ws.Range("E25").Value = "2026-02-17 22:19:32"
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H25").Value = "53%"
$ws.Range("I25").Copy()
$ws.Range("H25").PasteSpecial(-4122)
$ws.Range("E26").Value = "2026-02-17 22:19:35"
$ws.Range("E27").Value = "2026-02-17 22:19:37"
$ws.Range("E28").Value = "2026-02-17 22:19:39"
$ws.Range("J28").Value = "1018.5 hPa"
$ws.Range("O28").Value = "8.7 °C"
$ws.Range("E29").Value = "2026-02-17 22:19:42"
$ws.Range("E30").Value = "2026-02-17 22:19:44"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "66%"
$ws.Range("I30").Copy()
$ws.Range("H30").PasteSpecial(-4122)
$ws.Range("N30").Value = "7.7 °C 21:50 TU"
$ws.Range("O30").Value = "11.1 °C"
$ws.Range("E31").Value = "2026-02-17 22:19:46"
$ws.Range("E32").Value = "2026-02-17 22:19:49"
$ws.Range("E33").Value = "2026-02-17 22:19:51"
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = "45%"
$ws.Range("I33").Copy()
$ws.Range("H33").PasteSpecial(-4122)
$ws.Range("J33").Value = "1017.6 hPa"
$ws.Range("N33").Value = "3.0 °C 21:52 TU"
$ws.Range("E34").Value = "2026-02-17 22:19:54"
$ws.Range("N34").Value = "-1.9 °C 21:55 TU"
$ws.Range("O34").Value = "1.0 °C"
$ws.Range("E35").Value = "2026-02-17 22:19:56"
$ws.Range("J35").Value = "1020.2 hPa"
$ws.Range("E36").Value = "2026-02-17 22:19:59"
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = "61%"
$ws.Range("I36").Copy()
$ws.Range("H36").PasteSpecial(-4122)
$ws.Range("O36").Value = "12.3 °C"
$ws.Range("E37").Value = "2026-02-17 22:20:01"
$ws.Range("H37").NumberFormat = "@"
$ws.Range("H37").Value = "74%"
$ws.Range("I37").Copy()
$ws.Range("H37").PasteSpecial(-4122)
$ws.Range("N37").Value = "3.5 °C 21:59 TU"
$ws.Range("O37").Value = "7.2 °C"
$ws.Range("E38").Value = "2026-02-17 22:20:03"
$ws.Range("E39").Value = "2026-02-17 22:20:06"
$ws.Range("E40").Value = "2026-02-17 22:20:08"
$ws.Range("J40").Value = "1018.1 hPa"
$ws.Range("O40").Value = "9.4 °C"
$ws.Range("E41").Value = "2026-02-17 22:20:11"
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = "52%"
$ws.Range("I41").Copy()
$ws.Range("H41").PasteSpecial(-4122)
$ws.Range("E42").Value = "2026-02-17 22:20:13"
$ws.Range("H42").NumberFormat = "@"
$ws.Range("H42").Value = "62%"
$ws.Range("I42").Copy()
$ws.Range("H42").PasteSpecial(-4122)
$ws.Range("O42").Value = "12.6 °C"
$ws.Range("E43").Value = "2026-02-17 22:20:16"
$ws.Range("E44").Value = "2026-02-17 22:20:18"
$ws.Range("H44").NumberFormat = "@"
$ws.Range("H44").Value = "80%"
$ws.Range("I44").Copy()
$ws.Range("H44").PasteSpecial(-4122)
$ws.Range("M44").Value = "0.7 °C 21:55 TU"
$ws.Range("O44").Value = "-2.9 °C"
$ws.Range("E45").Value = "2026-02-17 22:20:20"
$ws.Range("H45").NumberFormat = "@"
$ws.Range("H45").Value = "72%"
$ws.Range("I45").Copy()
$ws.Range("H45").PasteSpecial(-4122)
$ws.Range("J45").Value = "1021.8 hPa"
$ws.Range("O45").Value = "5.2 °C"
$ws.Range("E46").Value = "2026-02-17 22:20:23"
$ws.Range("H46").NumberFormat = "@"
$ws.Range("H46").Value = "59%"
$ws.Range("I46").Copy()
$ws.Range("H46").PasteSpecial(-4122)
$ws.Range("N46").Value = "10.2 °C 21:39 TU"
$ws.Range("O46").Value = "15.1 °C"

$excel.CutCopyMode = $false
Write-Output "Applied 102 cell updates"
